$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 171 (shifts existing rows 171-280 down to 172-281)
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new price record
$ws.Range("A171").Value = 3
$ws.Range("B171").Value = "Femacal de La Calera"
$ws.Range("C171").Value = "Coquimbo"
$ws.Range("D171").Value = 45062
$ws.Range("E171").Value = 5
$ws.Range("F171").Value = 100112030
$ws.Range("G171").Value = "Poroto granado"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 70
$ws.Range("K171").Value = 31000
$ws.Range("L171").Value = 32000
$ws.Range("M171").Value = 31500
$ws.Range("N171").Value = "$/malla 25 kilos"
$ws.Range("O171").Value = "Provincia de Limarí"
$ws.Range("P171").Value = 1260
$ws.Range("Q171").Value = 25
$ws.Range("R171").Value = "Hortaliza"
